# Updated cryptos list on Sun Apr 28 14:59:22 UTC 2024 with GitHub Actions
# Refresh price/volume figures (and one rebranded coin row) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price strings must stay TEXT (not be reinterpreted as
# numbers by Excel), so we force a text number format before writing them
# and then restore the cell to its original (unstyled) "Normal" style.

$ws.Range("D2").Value = "63.634.13"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "3.314.88"
$ws.Range("E3").Value = "  +5.16%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "3.313.22"
$ws.Range("E8").Value = "  +5.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "3.855.66"
$ws.Range("E15").Value = "  +5.00%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "3.308.07"
$ws.Range("E17").Value = "  +5.01%  "
$ws.Range("D18").Value = "63.726.36"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.35%  "
$ws.Range("E24").Value = "  +5.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.43%  "
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.18%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "0.0₃0746"
$ws.Range("E38").Value = "  +6.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0402"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "432.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("D41").Value = "3.081.68"
$ws.Range("E41").Value = "  +5.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.265"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.89%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.115"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.03%  "
